$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows with only Price (D) and Volume (E) changes ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.029.86'
$ws.Range("E2").Value = '  -0.02%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.829.62'
$ws.Range("E3").Value = '  -0.16%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9980'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.43'
$ws.Range("E5").Value = '  +0.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6268'
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9994'
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07497'
$ws.Range("E8").Value = '  -1.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2918'
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.09'
$ws.Range("E10").Value = '  +2.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07705'
$ws.Range("E11").Value = '  -0.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.833.97'
$ws.Range("E12").Value = '  +0.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.988'
$ws.Range("E13").Value = '  +0.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6664'
$ws.Range("E14").Value = '  +0.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.43'
$ws.Range("E15").Value = '  -0.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009328'
$ws.Range("E16").Value = '  -7.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.979'
$ws.Range("E17").Value = '  -0.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.062.72'
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.080.12'
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.58'
$ws.Range("E20").Value = '  +1.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '223.12'
$ws.Range("E21").Value = '  -1.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.004'
$ws.Range("E22").Value = '  +0.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.102'
$ws.Range("E23").Value = '  -0.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9989'
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.46'
$ws.Range("E25").Value = '  +0.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1394'
$ws.Range("E26").Value = '  +1.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.486'
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.89'
$ws.Range("E28").Value = '  -0.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.496'
$ws.Range("E29").Value = '  +0.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05696'
$ws.Range("E30").Value = '  +9.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.144'
$ws.Range("E31").Value = '  +1.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.057'
$ws.Range("E32").Value = '  +1.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.200'
$ws.Range("E33").Value = '  +0.67%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7459'
$ws.Range("E34").Value = '  +1.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.840'
$ws.Range("E35").Value = '  -0.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.134'
$ws.Range("E36").Value = '  -0.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.665'
$ws.Range("E37").Value = '  -1.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.755'
$ws.Range("E38").Value = '  -0.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.217.08'
$ws.Range("E39").Value = '  -1.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.525'
$ws.Range("E41").Value = '  +3.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8899'
$ws.Range("E42").Value = '  -0.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9994'
$ws.Range("E43").Value = '  -0.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.98'
$ws.Range("E44").Value = '  +0.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.979.60'
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.040'
$ws.Range("E51").Value = '  +1.90%  '

# --- Rows with only Volume (E) changes ---
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("E50").Value = '  +1.04%  '

# --- Rows 46-49: coin entries reordered (Aave/BabyDogeCoin swap, Mantle/XinFinNetwork swap) ---
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.31'
$ws.Range("E46").Value = '  +1.88%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000123'
$ws.Range("E47").Value = '  -0.85%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5079'
$ws.Range("E48").Value = '  -0.55%  '
$ws.Range("B49").Value = 'XinFinNetwork'
$ws.Range("C49").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07565'
$ws.Range("E49").Value = '  +8.66%  '
